$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "[-, Joel L.-Tecnologia dos Materiais.-1A]"
$ws.Range("D2").Value = "[-, Valmir-Metrologia-1A, -, -]"
$ws.Range("F2").Value = "-"

# Row 3
$ws.Range("C3").Value = "[-, Valmir-Metrologia-1A, -, -]"
$ws.Range("D3").Value = "Cleidson-Circuitos elétricos"
$ws.Range("F3").Value = "-"

# Row 4
$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "[-, Valmir-Metrologia-1A, -, -]"
$ws.Range("D4").Value = "Cleidson-Circuitos elétricos"
$ws.Range("F4").Value = "-"

# Row 6
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = "[Joel L.-Tecnologia dos Materiais.-1A, André Guimarães-Desenho Técnico-1A]"
$ws.Range("D6").Value = "Cleidson-Circuitos elétricos"
$ws.Range("F6").Value = "-"

# Row 7
$ws.Range("B7").Value = "-"
$ws.Range("C7").Value = "Anselmo-Gestão integrada"
$ws.Range("D7").Value = "Cleidson-Circuitos elétricos"
$ws.Range("F7").Value = "-"

# Row 8
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "Anselmo-Gestão integrada"
$ws.Range("D8").Value = "[Valmir-Metrologia-1A, -, -, -]"
$ws.Range("F8").Value = "-"

# Row 11
$ws.Range("B11").Value = "-"
